$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row of data appended after the last existing row (row 66 -> new row 67).
# The Date column stores plain text like "MM/DD/YYYY" (not a real date), so a
# leading apostrophe forces text interpretation and avoids Excel's automatic
# date-value conversion; re-applying the "Normal" style afterwards keeps the
# cell on the workbook's default (unstyled) format, matching the other rows.
$ws.Range("A67").Value = "'01/30/2026"
$ws.Range("A67").Style = "Normal"

$ws.Range("B67").Value = 10956.4
$ws.Range("C67").Value = 0.237431419224339
$ws.Range("D67").Value = 0.762568580775661
$ws.Range("E67").Value = -236.89
$ws.Range("F67").Value = -30.18
$ws.Range("G67").Value = -22486.75
$ws.Range("H67").Value = -72.91
$ws.Range("I67").Value = -651.45
$ws.Range("J67").Value = -20.03
$ws.Range("K67").Value = -23138.2
$ws.Range("L67").Value = -67.86
